$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.905.96"
$ws.Range("E2").Value = "  +3.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.225.72"
$ws.Range("E3").Value = "  +3.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.11"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  +1.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.22"
$ws.Range("E7").Value = "  +1.65%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +2.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0882"
$ws.Range("E10").Value = "  +2.57%  "

$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.555.80"
$ws.Range("E12").Value = "  +3.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.14"
$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.44"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.827"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("E16").Value = "  +1.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.220.97"
$ws.Range("E17").Value = "  +2.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.805.10"
$ws.Range("E18").Value = "  +3.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.08"
$ws.Range("E19").Value = "  +3.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  +6.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.18"
$ws.Range("E21").Value = "  +0.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "256.84"
$ws.Range("E22").Value = "  +11.01%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -8.02%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "173.22"
$ws.Range("E26").Value = "  +0.67%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  +1.42%  "

$ws.Range("E28").Value = "  +5.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.38"
$ws.Range("E29").Value = "  +2.30%  "

$ws.Range("E30").Value = "  +2.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.85"
$ws.Range("E31").Value = "  +6.75%  "

$ws.Range("E32").Value = "  +1.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.68"
$ws.Range("E33").Value = "  +1.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("E34").Value = "  +2.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.24"
$ws.Range("E35").Value = "  +2.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0634"
$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("E37").Value = "  +6.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.49"
$ws.Range("E38").Value = "  +3.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.995"
$ws.Range("E39").Value = "  -0.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.90"
$ws.Range("E40").Value = "  +15.74%  "

$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.75"
$ws.Range("E42").Value = "  +12.25%  "

$ws.Range("B43").Value = "TerraClassic"
$ws.Range("C43").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000227"
$ws.Range("E43").Value = "  +52.89%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.24"
$ws.Range("E44").Value = "  -1.54%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.25"
$ws.Range("E45").Value = "  +5.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.63"
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.516.77"
$ws.Range("E47").Value = "  -1.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0946"
$ws.Range("E48").Value = "  +2.30%  "

$ws.Range("E49").Value = "  +2.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.84"
$ws.Range("E50").Value = "  +0.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.67"
$ws.Range("E51").Value = "  +11.71%  "

